$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -3
    7  = 4
    13 = -3
    17 = -6
    19 = -1
    21 = -4
    23 = 2
    24 = -2
    27 = -3
    30 = 0
    36 = -5
    46 = -6
    54 = -1
    56 = -4
    57 = -5
    65 = -7
    72 = -7
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
